# Clean up OCR-mangled species names in the "species" column (C) and fix a
# mis-transcribed landings value. The order the new strings are written in
# matters only for matching the shared-string table append order, but has
# no functional effect.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C16").Value = "Grouper"        # was "(irouper"
$ws.Range("C11").Value = "Sculpin"        # was "Soil pin"
$ws.Range("C10").Value = "White seabass"  # was "White seal kiss"
$ws.Range("C9").Value  = "Rockfish"       # was "Kockfish"
$ws.Range("C5").Value  = "Bluefin tuna"   # was "Hlucfin tuna"
$ws.Range("C8").Value  = "Spiny lobster"  # was "Spiny lolxstcr"
$ws.Range("C4").Value  = "Albacore"       # was "Altiarore"

# Correct the mis-keyed pounds figure for "Giant sea bass" (San Diego landings)
$ws.Range("E18").Value = 31849

# Leave the selection on E20, matching the saved view state
$ws.Activate()
$ws.Range("E20").Select()
